# Update cryptocurrency price/volume data per the data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.937.49"
$ws.Range("E2").Value = "  +6.47%  "

$ws.Range("D3").Value = "3.015.23"
$ws.Range("E3").Value = "  +3.79%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'584.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").Value = "'161.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.42%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.011.95"
$ws.Range("E8").Value = "  +3.77%  "

$ws.Range("E9").Value = "  +3.88%  "

$ws.Range("D10").Value = "'6.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("E11").Value = "  +6.08%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.93%  "

$ws.Range("E13").Value = "  +9.02%  "

$ws.Range("D14").Value = "'34.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.71%  "

$ws.Range("D16").Value = "65.923.50"
$ws.Range("E16").Value = "  +6.59%  "

$ws.Range("D17").Value = "3.515.04"
$ws.Range("E17").Value = "  +3.77%  "

$ws.Range("D18").Value = "'6.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.01%  "

$ws.Range("D19").Value = "3.016.08"
$ws.Range("E19").Value = "  +4.52%  "

$ws.Range("D20").Value = "'456.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.47%  "

$ws.Range("E21").Value = "  +7.33%  "

$ws.Range("D22").Value = "'0.688"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.74%  "

$ws.Range("E23").Value = "  +7.67%  "

$ws.Range("D24").Value = "'82.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.65%  "

$ws.Range("E25").Value = "  +12.44%  "

$ws.Range("D26").Value = "'12.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.17%  "

$ws.Range("D27").Value = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.79%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "'8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.41%  "

$ws.Range("E30").Value = "  +16.91%  "

$ws.Range("E31").Value = "  -6.53%  "

$ws.Range("E32").Value = "  +4.10%  "

$ws.Range("D33").Value = "'27.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.29%  "

$ws.Range("D34").Value = "'0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").Value = "'0.994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.10%  "

$ws.Range("E37").Value = "  +8.64%  "

$ws.Range("D38").Value = "'2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.45%  "

$ws.Range("D39").Value = "'49.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("E40").Value = "  +2.64%  "

$ws.Range("D41").Value = "'0.311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.65%  "

$ws.Range("E42").Value = "  +7.11%  "

$ws.Range("D43").Value = "'43.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.74%  "

$ws.Range("E44").Value = "  +3.62%  "

$ws.Range("D45").Value = "'391.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.91%  "

$ws.Range("D46").Value = "2.798.93"
$ws.Range("E46").Value = "  +3.42%  "

$ws.Range("E47").Value = "  +5.64%  "

$ws.Range("D48").Value = "'134.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("D50").Value = "'23.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.31%  "

$ws.Range("E51").Value = "  +4.53%  "

